$wb = $excel.ActiveWorkbook

# Rename sheets: Foglio1 -> Data, Foglio2 -> Legend
$data = $wb.Worksheets.Item("Foglio1")
$data.Name = "Data"
$legend = $wb.Worksheets.Item("Foglio2")
$legend.Name = "Legend"

# Populate the "Legend" sheet with a Field / Description glossary table
# mirroring the header row of the "Data" sheet.
$legend.Range("A1").Value = "Field"
$legend.Range("B1").Value = "Description"

$legend.Range("A2").Value = "IP/Domain"
$legend.Range("B2").Value = "The IP or Domain analyzed"

$legend.Range("A3").Value = "URL Domain Report"
$legend.Range("B3").Value = "The VirusTotal ""Relations"" Tab of the analyzed Domain"

$legend.Range("A4").Value = "Name or hash of detection"
$legend.Range("B4").Value = "The ""Communicating file"" found in ""Relations"" Tab of VirusTotal, for the analyzed domain"

$legend.Range("A5").Value = "URL hash analysis"
$legend.Range("B5").Value = "The VirusTotal scan page for the analyzed communicating file"

$legend.Range("A6").Value = "Scanned Data"
$legend.Range("B6").Value = "The data reported on ""Communicating Files"" table, in the ""Relations"" Tab, for the analyzed domain"

$legend.Range("A7").Value = "Detection"
$legend.Range("B7").Value = "Number of engines that found the communicating file malicious"

$legend.Range("A8").Value = "Total of scanners"
$legend.Range("B8").Value = "Total engines reported on the result page on VirusTotal, for the communicating file analyzed"

$legend.Range("A9").Value = "Malicious contacted IP"
$legend.Range("B9").Value = "Record of ""Contacted IP"" table, in the ""Relations"" Tab, on Communicating file analysis"

$legend.Range("A10").Value = "Is this IP malicious?"
$legend.Range("B10").Value = "If there are founds on the ""Reverse.it"" portal, reporting malicious files, then set to ""yes"""

$legend.Range("A11").Value = "URL Reverse Report (require an Account"
$legend.Range("B11").Value = "URL of Reverse.it portal, with the evidence the IP is malicious"

# Yellow highlight on the "Field" column, matching the look of the header
# row on the "Data" sheet (same fill used there).
$legend.Range("A2:A6").Interior.Color = 65535
$legend.Range("A7:A11").Interior.Color = 65535

# Column widths sized to fit the longest content in each column.
$legend.Columns.Item(1).ColumnWidth = 36.15
$legend.Columns.Item(2).ColumnWidth = 89.65

# Leave the selection parked below the table, like the source workbook,
# then bring the focus back to "Data" so it stays the active tab.
$legend.Range("A12").Select()
$data.Activate()
